$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 06:22"

# Update Arabia Saudita row (row 39) stats
$ws.Range("B39").Value = 2370
$ws.Range("C39").Value = 191
$ws.Range("E39").Value = 1921

# Row 152 now holds "San Martin (Parte Holandesa)" with fresh data
$ws.Range("A152").Value = "San Martin (Parte Holandesa)"
$ws.Range("B152").Value = 25
$ws.Range("C152").Value = 2
$ws.Range("D152").Value = 6
$ws.Range("E152").Value = 15

# Row 153 now holds "Guyana" (previous Guyana figures shifted down)
$ws.Range("A153").Value = "Guyana"
$ws.Range("C153").Value = 1
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 20
$ws.Range("F153").Value = 0
$ws.Range("H153").Value = 4

# Row 154 now holds "San Martin (Parte Francesa)" (previous figures shifted down)
$ws.Range("A154").Value = "San Martin (Parte Francesa)"
$ws.Range("B154").Value = 24
$ws.Range("D154").Value = 5
$ws.Range("E154").Value = 17
$ws.Range("F154").Value = 6
$ws.Range("H154").Value = 2
